$d = $word.ActiveDocument

# The document currently has two consecutive empty paragraphs right
# before the "Version 0.1" heading. The second of those two empty
# paragraphs is where the new "Version 0.2" changelog section needs to
# go (the first empty paragraph is left untouched as a spacer).
#
# We locate it by walking the paragraphs and looking for the second
# empty paragraph that immediately precedes the paragraph whose text
# starts with "Version 0.1".
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "Version 0.1" -and $i -ge 2) {
        $prev = $d.Paragraphs.Item($i - 1)
        if ($prev.Range.Text.Trim() -eq "") {
            $targetIndex = $i - 1
        }
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the empty paragraph preceding 'Version 0.1'"
}

$target = $d.Paragraphs.Item($targetIndex)

# Build the replacement OOXML: a Heading1 "Version 0.2" title, the
# revision date, and three bulleted changelog entries using the same
# list (numId 23) as the existing "Version 0.1" bullet.
$xml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Version 0.</w:t></w:r><w:r><w:t>2</w:t></w:r></w:p><w:p><w:r><w:t>12-11-2021</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="23"/></w:numPr></w:pPr><w:r><w:t>Shortened the width of the side parts with 0.2mm for a better fit</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="23"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Reduced the width of the </w:t></w:r><w:r><w:t>LCD gap in the front for a better fit</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="23"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Reduced the holes for the screws holding the </w:t></w:r><w:r><w:t>SDS011 to fit 3mm screws</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# InsertXML replaces the contents of the target range with the
# supplied OOXML fragment, turning the single empty paragraph into the
# five new paragraphs above.
[void]$target.Range.InsertXML($xml)
